# Update countries & provincias Spain
# The underlying COVID data feed refreshed (00:40 -> 01:10). Several countries'
# totals grew enough to change their rank in the table (which is kept sorted
# descending by "Casos totales" / column B), so besides updating individual
# countries' figures we also need to re-seat the rows that swapped order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $nuevos
    $ws.Cells.Item($row, 4).Value = $activos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $criticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# Header: data refresh timestamp
$ws.Range("A1").Value = "Datos actualizados a 29 de Mayo de 2020 a las 01:10"

# Estados Unidos keeps rank #1 but the figures moved on
Set-Row 4 "Estados Unidos" 1767577 21774 498646 1165629 0 1195 103302

# Colombia overtook Irlanda, Indonesia and Kuwait, so the block shifts down
# one row and Colombia's refreshed figures take the top spot
Set-Row 34 "Colombia" 25366 1262 6665 17879 0 19 822
Set-Row 35 "Irlanda" 24841 38 22089 1113 0 8 1639
Set-Row 36 "Indonesia" 24538 687 6240 16802 0 23 1496
Set-Row 37 "Kuwait" 24112 845 8698 15229 0 10 185

# In-place figure refreshes (rank unchanged)
Set-Row 47 "Argentina" 14702 769 4617 9577 0 8 508
Set-Row 55 "Chequia" 9140 54 6460 2361 0 2 319
Set-Row 58 "Nigeria" 8915 182 2592 6064 0 5 259

# Venezuela overtook Haiti
Set-Row 104 "Venezuela" 1325 80 302 1012 0 0 11
Set-Row 105 "Haiti" 1320 146 22 1264 0 1 34

# Tied totals re-ordered amongst themselves (Fiyi/Curazao, Belice/Santa Lucia)
Set-Row 198 "Fiyi" 18 0 15 3 0 0 0
Set-Row 199 "Curazao" 18 0 14 3 0 0 1
Set-Row 200 "Belice" 18 0 16 0 0 0 2
Set-Row 201 "Santa Lucia" 18 0 18 0 0 0 0

# Montserrat/Seychelles swap
Set-Row 210 "Montserrat" 11 0 10 0 0 0 1
Set-Row 211 "Seychelles" 11 0 11 0 0 0 0

# San Bartolome/Bonaire swap
Set-Row 215 "San Bartolome" 6 0 6 0 0 0 0
Set-Row 216 "Bonaire, San Eustaquio y Saba" 6 0 6 0 0 0 0
